$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add a new student worker ("Ajay") on row 13 ---
$ws.Range("A13").Value = "Ajay"
$ws.Range("B13").Value = 15

# Mon Start / Mon End  (9:00 AM - 12:00 PM)
$ws.Range("C13").Value = 0.375
$ws.Range("D13").Value = 0.5

# Tue Start / Tue End  (9:00 AM - 12:00 PM)
$ws.Range("K13").Value = 0.375
$ws.Range("L13").Value = 0.5

# Wed Start / Wed End  (8:30 AM - 11:30 AM)
$ws.Range("O13").Value = 0.35416666666666669
$ws.Range("P13").Value = 0.47916666666666669

# Thur Start / Thur End (9:00 AM - 3:30 PM)
$ws.Range("S13").Value = 0.375
$ws.Range("T13").Value = 0.64583333333333337

# The "start" columns for this row should pick up the themed font colour
# (matching the other rows) instead of the default explicit-black font;
# restoring the cell style to "Normal" and re-applying the time format
# gets them back onto the same (themed-font) style used elsewhere.
$ws.Range("C13").Style = "Normal"
$ws.Range("C13").NumberFormat = "h:mm AM/PM"
$ws.Range("D13").Style = "Normal"
$ws.Range("D13").NumberFormat = "h:mm AM/PM"
$ws.Range("K13").Style = "Normal"
$ws.Range("K13").NumberFormat = "h:mm AM/PM"
$ws.Range("L13").Style = "Normal"
$ws.Range("L13").NumberFormat = "h:mm AM/PM"
$ws.Range("O13").Style = "Normal"
$ws.Range("O13").NumberFormat = "h:mm AM/PM"
$ws.Range("S13").Style = "Normal"
$ws.Range("S13").NumberFormat = "h:mm AM/PM"

# --- Re-freeze the panes one column to the left and move the selection ---
$ws.Range("B1").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D4").Select() | Out-Null
